$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.046.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.631.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.052.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0749"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("E22").Value = "  -4.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0517"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "

$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.791"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +44.19%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("E33").Value = "  -0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.337.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.22%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.755.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.903"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +34.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.69%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
